$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.537.55"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "1.474.52"
$ws.Range("E3").Value = "  +2.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9543"
$ws.Range("E5").Value = "  +3.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.77"
$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3616"
$ws.Range("E7").Value = "  -1.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3064"
$ws.Range("E8").Value = "  -1.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.46"
$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.063"
$ws.Range("E10").Value = "  +3.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06654"
$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.526"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.15"
$ws.Range("E14").Value = "  +2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.192"
$ws.Range("E15").Value = "  +1.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9539"
$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001029"
$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").Value = "1.474.83"
$ws.Range("E18").Value = "  +1.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05931"
$ws.Range("E19").Value = "  +5.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.23"
$ws.Range("E20").Value = "  +2.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.502"
$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.50"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.17"
$ws.Range("E23").Value = "  +2.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.258"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").Value = "20.569.08"
$ws.Range("E25").Value = "  +1.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.15"
$ws.Range("E26").Value = "  +5.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.142"
$ws.Range("E27").Value = "  -2.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.20"
$ws.Range("E28").Value = "  +0.95%  "

$ws.Range("D29").Value = "1.636.02"
$ws.Range("E29").Value = "  +2.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.72"
$ws.Range("E30").Value = "  +2.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.937"
$ws.Range("E31").Value = "  +4.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.989"
$ws.Range("E32").Value = "  +2.97%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8076"
$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07990"
$ws.Range("E34").Value = "  +4.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.513"
$ws.Range("E35").Value = "  +1.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.222"
$ws.Range("E36").Value = "  +8.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05847"
$ws.Range("E37").Value = "  -3.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.727"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02054"
$ws.Range("E39").Value = "  +3.13%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.37"
$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9549"
$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1879"
$ws.Range("E42").Value = "  +2.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.418"
$ws.Range("E43").Value = "  +4.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5293"
$ws.Range("E44").Value = "  +0.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.524"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.22"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.38"
$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5203"
$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.816"
$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9791"
$ws.Range("E51").Value = "  -1.49%  "
